$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (rows 2..97) holds labels "q1".."q96".
# Renumber each label down by one: "q1" -> "q0", "q2" -> "q1", ..., "q96" -> "q95".
for ($r = 2; $r -le 97; $r++) {
    $n = $r - 2
    $ws.Cells.Item($r, 1).Value = "q$n"
}
